$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add a new "2022-Q1" sheet, positioned between "2021-Q4" and "总计".
#    Copying "2021-Q4" gives us the same layout/formatting (header
#    style, column widths, page margins, etc.) as the other quarterly
#    sheets, so we only need to overwrite the cell values afterwards.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Header row (unchanged labels, just rewritten to be explicit)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Make sure the fund-code / numeric-looking text columns stay text
# (so leading zeros and exact decimal strings are preserved).
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "008763"
$newSheet.Range("C2").Value = "天弘越南市场股票（QDII）A"
$newSheet.Range("D2").Value = "37.53"
$newSheet.Range("E2").Value = "92.10"
$newSheet.Range("F2").Value = "5.93"
$newSheet.Range("G2").Value = "2.2255"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "008764"
$newSheet.Range("C3").Value = "天弘越南市场股票（QDII）C"
$newSheet.Range("D3").Value = "14.26"
$newSheet.Range("E3").Value = "92.10"
$newSheet.Range("F3").Value = "5.93"
$newSheet.Range("G3").Value = "0.8456"
$newSheet.Range("H3").Value = 4

# ------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the
#    existing rows down and bumping their running index in column A.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Copy formatting from the row that just got pushed down to row 3.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 3.07

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the originally-active sheet/tab (copying "2021-Q4" above made
# the new sheet active as a side effect).
$wb.Worksheets.Item("2021-Q2").Activate()
